$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14 mirrors the pattern of the other "Board/Function" rows above it.
# Shared strings must be introduced in the same order as the target file:
# "none" (R14/S14) first, then "AltFunc" (B14), so they land at sharedString
# indices 29 and 30 respectively.
$ws.Range("R14").Value = "none"
$ws.Range("S14").Value = "none"
$ws.Range("B14").Value = "AltFunc"

$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 13
$ws.Range("E14").Value = 14
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = 16
$ws.Range("H14").Value = 17
$ws.Range("I14").Value = 18
$ws.Range("J14").Value = 19
$ws.Range("K14").Value = 20
$ws.Range("L14").Value = 21
$ws.Range("M14").Value = 22

# Update selection to match the new active cell.
$ws.Range("B14").Select()
